# Applies numeric corrections captured in the commit diff for Garuda_Profits.xlsx
# (market-price / profit recalculations per leve row across multiple job sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3991.8667
$ws.Range("I41").Value = 91.28570999999999
$ws.Range("J41").Value = 7404.875
$ws.Range("K41").Value = 91.28570999999999
$ws.Range("L41").Value = 7404.875
$ws.Range("M41").Value = 348.71429
$ws.Range("N41").Value = -8284.875
$ws.Range("H69").Value = 3452.258
$ws.Range("I69").Value = 3454
$ws.Range("J69").Value = 3400
$ws.Range("K69").Value = 10362
$ws.Range("L69").Value = 10200
$ws.Range("M69").Value = -9488
$ws.Range("N69").Value = -11948
$ws.Range("H72").Value = 3452.258
$ws.Range("I72").Value = 3454
$ws.Range("J72").Value = 3400
$ws.Range("K72").Value = 31086
$ws.Range("L72").Value = 30600
$ws.Range("M72").Value = -26718
$ws.Range("N72").Value = -39336
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 833.3333
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 833.3333
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -292.3333
$ws.Range("N100").Value = -3082
$ws.Range("H116").Value = 1915.3846
$ws.Range("I116").Value = 1922.2222
$ws.Range("K116").Value = 1922.2222
$ws.Range("M116").Value = 1519.7778
$ws.Range("H129").Value = 29710.057
$ws.Range("I129").Value = 560.26666
$ws.Range("J129").Value = 51572.4
$ws.Range("K129").Value = 1680.79998
$ws.Range("L129").Value = 154717.2
$ws.Range("M129").Value = 3319.20002
$ws.Range("N129").Value = -164717.2
$ws.Range("H135").Value = 625.8333
$ws.Range("I135").Value = 522.3111
$ws.Range("J135").Value = 2178.6667
$ws.Range("K135").Value = 4700.7999
$ws.Range("L135").Value = 19608.0003
$ws.Range("M135").Value = -2165.7999
$ws.Range("N135").Value = -24678.0003
$ws.Range("H141").Value = 1544.7424
$ws.Range("I141").Value = 930.5349
$ws.Range("J141").Value = 2693.0435
$ws.Range("K141").Value = 2791.6047
$ws.Range("L141").Value = 8079.130500000001
$ws.Range("M141").Value = 2388.3953
$ws.Range("N141").Value = -18439.1305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 212.45454
$ws.Range("I80").Value = 540
$ws.Range("J80").Value = 160.73685
$ws.Range("K80").Value = 540
$ws.Range("L80").Value = 160.73685
$ws.Range("M80").Value = 458
$ws.Range("N80").Value = -2156.73685
$ws.Range("H82").Value = 2511.6667
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 212.45454
$ws.Range("I83").Value = 540
$ws.Range("J83").Value = 160.73685
$ws.Range("K83").Value = 2700
$ws.Range("L83").Value = 803.68425
$ws.Range("M83").Value = 2292
$ws.Range("N83").Value = -10787.68425
$ws.Range("H85").Value = 2511.6667
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 4208.4287
$ws.Range("I86").Value = 2927.25
$ws.Range("J86").Value = 5916.6665
$ws.Range("K86").Value = 2927.25
$ws.Range("L86").Value = 5916.6665
$ws.Range("M86").Value = -1804.25
$ws.Range("N86").Value = -8162.6665
$ws.Range("H89").Value = 4208.4287
$ws.Range("I89").Value = 2927.25
$ws.Range("J89").Value = 5916.6665
$ws.Range("K89").Value = 14636.25
$ws.Range("L89").Value = 29583.3325
$ws.Range("M89").Value = -9020.25
$ws.Range("N89").Value = -40815.3325
$ws.Range("H134").Value = 24603.045
$ws.Range("I134").Value = 28614.703
$ws.Range("J134").Value = 3398.5715
$ws.Range("K134").Value = 85844.109
$ws.Range("L134").Value = 10195.7145
$ws.Range("M134").Value = -83309.109
$ws.Range("N134").Value = -15265.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5052961.5
$ws.Range("I31").Value = 1664.72
$ws.Range("K31").Value = 1664.72
$ws.Range("M31").Value = -1369.72
$ws.Range("H34").Value = 5052961.5
$ws.Range("I34").Value = 1664.72
$ws.Range("K34").Value = 1664.72
$ws.Range("M34").Value = -1462.72
$ws.Range("H58").Value = 764.57574
$ws.Range("I58").Value = 896.1111
$ws.Range("J58").Value = 606.73334
$ws.Range("K58").Value = 896.1111
$ws.Range("L58").Value = 606.73334
$ws.Range("M58").Value = -693.1111
$ws.Range("N58").Value = -1012.73334
$ws.Range("H74").Value = 29864.6
$ws.Range("J74").Value = 29864.6
$ws.Range("L74").Value = 29864.6
$ws.Range("N74").Value = -31612.6
$ws.Range("H77").Value = 29864.6
$ws.Range("J77").Value = 29864.6
$ws.Range("L77").Value = 89593.79999999999
$ws.Range("N77").Value = -98329.79999999999
$ws.Range("H132").Value = 1357.0266
$ws.Range("I132").Value = 1277.9636
$ws.Range("K132").Value = 3833.8908
$ws.Range("M132").Value = -1303.8908
$ws.Range("H136").Value = 764.57574
$ws.Range("I136").Value = 896.1111
$ws.Range("J136").Value = 606.73334
$ws.Range("K136").Value = 2688.3333
$ws.Range("L136").Value = 1820.20002
$ws.Range("M136").Value = -138.3332999999998
$ws.Range("N136").Value = -6920.20002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5781.5713
$ws.Range("I131").Value = 7336.25
$ws.Range("J131").Value = 806.6
$ws.Range("K131").Value = 22008.75
$ws.Range("L131").Value = 2419.8
$ws.Range("M131").Value = -16968.75
$ws.Range("N131").Value = -12499.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3837.6365
$ws.Range("I80").Value = 2536.8572
$ws.Range("J80").Value = 6114
$ws.Range("K80").Value = 2536.8572
$ws.Range("L80").Value = 6114
$ws.Range("M80").Value = -1538.8572
$ws.Range("N80").Value = -8110
$ws.Range("H83").Value = 3837.6365
$ws.Range("I83").Value = 2536.8572
$ws.Range("J83").Value = 6114
$ws.Range("K83").Value = 12684.286
$ws.Range("L83").Value = 30570
$ws.Range("M83").Value = -7692.286
$ws.Range("N83").Value = -40554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 267.6087
$ws.Range("I55").Value = 266.66666
$ws.Range("K55").Value = 266.66666
$ws.Range("M55").Value = -93.66665999999998
$ws.Range("H93").Value = 2253903
$ws.Range("I93").Value = 2458621
$ws.Range("J93").Value = 2004
$ws.Range("K93").Value = 2458621
$ws.Range("L93").Value = 2004
$ws.Range("M93").Value = -2457373
$ws.Range("N93").Value = -4500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10280

